$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.630.09'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '3.509.06'
$ws.Range("E3").Value = '  -1.99%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.60'
$ws.Range("E5").Value = '  -3.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.18'
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("D7").Value = '3.508.98'
$ws.Range("E7").Value = '  -1.99%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.27'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '4.109.28'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.65'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("D17").Value = '3.512.56'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '64.607.83'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.19'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("E21").Value = '  -3.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.20'
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.578'
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("D24").Value = '3.651.56'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.98'
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000110'
$ws.Range("E27").Value = '  -4.97%  '
$ws.Range("E28").Value = '  -5.02%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -7.64%  '
$ws.Range("E31").Value = '  -5.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -4.75%  '
$ws.Range("D33").Value = '3.514.33'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '24.01'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.27'
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '171.00'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.95'
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("E42").Value = '  -2.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.46'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("E45").Value = '  -2.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.39'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D49").Value = '2.459.10'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.88'
$ws.Range("E50").Value = '  -1.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.893'
$ws.Range("E51").Value = '  +2.01%  '
